# B1--and-B2-PowerPoint.pptx edit
#
# 1) Slide 5's table (graphicFrame "Google Shape;122;p17") switches its
#    table style from {E9224E8E-D522-470A-A2D1-EFA65EFD3AD3} to
#    {A6A97345-076E-4D6C-ABD9-6A59C38B2474}.
# 2) The deck's theme reverts from the "Integral" palette to the classic
#    "Office Theme" palette (the font scheme and format scheme - fills,
#    lines, effects - are already identical between the two themes, only
#    the 12 theme colours differ).

$p = $ppt.ActivePresentation

# --- 1) table style -------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{A6A97345-076E-4D6C-ABD9-6A59C38B2474}")

# --- 2) theme colours -------------------------------------------------------
# PowerPoint's ThemeColor.RGB uses the OLE/COM BGR-packed long (0x00BBGGRR),
# so the familiar RRGGBB hex of the "Office Theme" palette needs swapping
# to bbggrr before it is assigned.
$officeThemeColors = @{
    1  = 0x000000   # dk1      000000
    2  = 0xFFFFFF   # lt1      FFFFFF
    3  = 0x6A5444   # dk2      44546A
    4  = 0xE6E6E7   # lt2      E7E6E6
    5  = 0xD59B5B   # accent1  5B9BD5
    6  = 0x317DED   # accent2  ED7D31
    7  = 0xA5A5A5   # accent3  A5A5A5
    8  = 0x00C0FF   # accent4  FFC000
    9  = 0xC47244   # accent5  4472C4
    10 = 0x47AD70   # accent6  70AD47
    11 = 0xC16305   # hlink    0563C1
    12 = 0x724F95   # folHlink 954F72
}

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i]
}
